$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the old (empty) row 83, pushing everything
# below down by two rows. The new rows pick up row 82's formatting.
$ws.Range("A83:D84").EntireRow.Insert()

# Copy row 82's cell formatting (date style in col A, plain style in B:D)
# down into the two freshly inserted rows so the new date cells render
# with the same number format as the rest of the date column.
$ws.Range("A82:D82").Copy()
$ws.Range("A83:D83").PasteSpecial(-4122)
$ws.Range("A84:D84").PasteSpecial(-4122)
$ws.Range("A82").Copy()
$ws.Range("A84").PasteSpecial(-4122)

# New match results for 2018-11-21 (serial 43417): Fritz 15-13 on side A,
# then the rematch 15-13 on side H.
$ws.Range("A83").Value = 43417
$ws.Range("B83").Value = 15
$ws.Range("C83").Value = 13
$ws.Range("D83").Value = "A"

$ws.Range("A84").Value = 43417
$ws.Range("B84").Value = 15
$ws.Range("C84").Value = 13
$ws.Range("D84").Value = "H"

# Column A now gets a bit wider to fit the date values, with an explicit
# custom width separate from the rest of the table columns.
$ws.Columns("A").ColumnWidth = 8.67

# Restore the selection to where editing left off (no more pinned
# top-left scroll position).
$ws.Range("H88").Select()
